$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AI3").Value = "Passed"
$ws.Range("AI4").Value = "Passed"
$ws.Range("AI5").Value = "Passed"
$ws.Range("AI6").Value = "Passed"
